# Add new loan records (rows 23-27) to the historial_prestamos sheet,
# matching the rows already present (EQ equipment loans for Jhonathan / Oficial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the date/time style already applied to the existing "Hora_Prestamo"
# column (G) so new timestamp cells get the same number format.
$dateFormat = $ws.Cells.Item(22, 7).NumberFormat

$newRows = @(
    @{ Row=23; A=12345; B="Jhonathan"; C="Oficial"; D="EQ-002"; E="Laptop Dell";        G=45776.66991229167; I="Prestado"; J="nan" },
    @{ Row=24; A=12345; B="Jhonathan"; C="Oficial"; D="EQ-010"; E="Radio Comunicador";  G=45776.66991229167; I="Prestado"; J="nan" },
    @{ Row=25; A=12345; B="Jhonathan"; C="Oficial"; D="EQ-001"; E="Laptop Dell";        G=45776.68377900463; I="Prestado"; J="nan" },
    @{ Row=26; A=12345; B="Jhonathan"; C="Oficial"; D="EQ-009"; E="Radio Comunicador";  G=45776.68377900463; I="Prestado"; J="nan" },
    @{ Row=27; A=12345; B="Jhonathan"; C="Oficial"; D="EQ-007"; E="Proyector Epson";    G=45791.67583650942; I="Prestado"; J="" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 7).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 9).Value = $r.I
    if ($r.J -ne "") {
        $ws.Cells.Item($row, 10).Value = $r.J
    }
}
